$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "61.047.70"
$ws.Range("E2").Value = "  -3.80%  "

$ws.Range("D3").Value = "2.967.52"
$ws.Range("E3").Value = "  -3.75%  "

$ws.Range("E4").Value = "  +0.18%  "

Set-TextValue "D5" "542.86"
$ws.Range("E5").Value = "  -0.50%  "

Set-TextValue "D6" "128.44"
$ws.Range("E6").Value = "  -7.99%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "2.962.07"
$ws.Range("E8").Value = "  -3.78%  "

Set-TextValue "D9" "0.489"
$ws.Range("E9").Value = "  -2.51%  "

$ws.Range("E10").Value = "  -7.13%  "

Set-TextValue "D11" "0.142"
$ws.Range("E11").Value = "  -8.66%  "

Set-TextValue "D12" "0.437"
$ws.Range("E12").Value = "  -4.65%  "

Set-TextValue "D13" "0.0000215"
$ws.Range("E13").Value = "  -4.15%  "

Set-TextValue "D14" "33.24"
$ws.Range("E14").Value = "  -4.99%  "

$ws.Range("D15").Value = "3.463.04"
$ws.Range("E15").Value = "  -3.41%  "

$ws.Range("D16").Value = "61.177.23"
$ws.Range("E16").Value = "  -3.63%  "

$ws.Range("E17").Value = "  -2.99%  "

$ws.Range("D18").Value = "2.976.95"
$ws.Range("E18").Value = "  -3.49%  "

Set-TextValue "D19" "6.50"
$ws.Range("E19").Value = "  -2.41%  "

Set-TextValue "D20" "467.48"
$ws.Range("E20").Value = "  -1.55%  "

Set-TextValue "D21" "12.93"
$ws.Range("E21").Value = "  -4.21%  "

Set-TextValue "D22" "0.654"
$ws.Range("E22").Value = "  -6.91%  "

Set-TextValue "D23" "6.85"
$ws.Range("E23").Value = "  -3.45%  "

Set-TextValue "D24" "79.10"
$ws.Range("E24").Value = "  +0.44%  "

Set-TextValue "D25" "11.81"
$ws.Range("E25").Value = "  -3.65%  "

Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("E27").Value = "  -1.80%  "

Set-TextValue "D28" "7.51"
$ws.Range("E28").Value = "  -5.90%  "

$ws.Range("E29").Value = "  -0.06%  "

Set-TextValue "D30" "1.86"
$ws.Range("E30").Value = "  -2.45%  "

Set-TextValue "D31" "25.18"
$ws.Range("E31").Value = "  -4.00%  "

$ws.Range("E32").Value = "  -3.72%  "

Set-TextValue "D33" "2.26"
$ws.Range("E33").Value = "  -2.96%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D34" "54.18"
$ws.Range("E34").Value = "  -6.73%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D35" "5.32"
$ws.Range("E35").Value = "  -1.89%  "

Set-TextValue "D36" "5.77"
$ws.Range("E36").Value = "  -4.08%  "

Set-TextValue "D37" "446.37"
$ws.Range("E37").Value = "  -9.69%  "

$ws.Range("D38").Value = "3.101.09"
$ws.Range("E38").Value = "  -4.94%  "

Set-TextValue "D39" "0.0778"
$ws.Range("E39").Value = "  -2.70%  "

Set-TextValue "D40" "0.0372"
$ws.Range("E40").Value = "  -7.70%  "

Set-TextValue "D41" "0.114"
$ws.Range("E41").Value = "  -2.86%  "

Set-TextValue "D42" "7.94"
$ws.Range("E42").Value = "  -2.57%  "

$ws.Range("E43").Value = "  -0.05%  "

Set-TextValue "D44" "2.26"
$ws.Range("E44").Value = "  -13.53%  "

Set-TextValue "D45" "24.98"
$ws.Range("E45").Value = "  -1.59%  "

Set-TextValue "D46" "0.236"
$ws.Range("E46").Value = "  -7.44%  "

Set-TextValue "D47" "0.106"
$ws.Range("E47").Value = "  -3.04%  "

$ws.Range("E48").Value = "  +10.17%  "

Set-TextValue "D49" "1.89"
$ws.Range("E49").Value = "  -7.09%  "

Set-TextValue "D50" "113.38"
$ws.Range("E50").Value = "  -8.58%  "

$ws.Range("D51").Value = "0.0₃0469"
$ws.Range("E51").Value = "  -11.49%  "
